{"js": "// Replace the date line and every \"a OP b=\" arithmetic-problem cell in the\n// single table with the values from the target revision. Replacements are\n// applied strictly by document position (first paragraph, then table cells\n// in row-major order) since several cell strings are not unique (e.g.\n// \"85-28=\" and \"60+11=\" each occur twice in the source), so plain\n// find/replace-by-text would be ambiguous.\n//\n// [oldText, newText] pairs in document order: index 0 is the date\n// paragraph, indices 1..100 are the 20x5 table cells read left-to-right,\n// top-to-bottom.\nconst replacements = [[\"2023-08-12 Saturday\", \"2023-08-13 Sunday\"], [\"13+0=\", \"0+26=\"], [\"16-11=\", \"93-65=\"], [\"93-4=\", \"60-12=\"], [\"26+9=\", \"92-43=\"], [\"93-71=\", \"74-22=\"], [\"82+2=\", \"74+12=\"], [\"58+3=\", \"21+75=\"], [\"50-14=\", \"0+40=\"], [\"66+17=\", \"11+31=\"], [\"13+18=\", \"76-46=\"], [\"9+40=\", \"68+31=\"], [\"75+8=\", \"66+8=\"], [\"83+1=\", \"18+20=\"], [\"38-29=\", \"14+81=\"], [\"56-15=\", \"32+24=\"], [\"0+90=\", \"52+18=\"], [\"11+66=\", \"37+37=\"], [\"87-19=\", \"52-49=\"], [\"80-17=\", \"70+25=\"], [\"20+34=\", \"30+47=\"], [\"77+12=\", \"47-29=\"], [\"57-37=\", \"9+75=\"], [\"48+48=\", \"53+22=\"], [\"12+79=\", \"81-34=\"], [\"67-42=\", \"16+12=\"], [\"72-20=\", \"67-22=\"], [\"27-22=\", \"19+8=\"], [\"88-40=\", \"1+53=\"], [\"8+70=\", \"23+18=\"], [\"46-24=\", \"30+68=\"], [\"96-67=\", \"31+64=\"], [\"91-46=\", \"88-29=\"], [\"41+58=\", \"68-59=\"], [\"2+22=\", \"39+11=\"], [\"89+3=\", \"9+25=\"], [\"0+99=\", \"99-58=\"], [\"67-1=\", \"99-39=\"], [\"26+7=\", \"2+66=\"], [\"96+2=\", \"93+3=\"], [\"16-7=\", \"65-27=\"], [\"67-40=\", \"53-37=\"], [\"69-39=\", \"83-5=\"], [\"89-60=\", \"11-10=\"], [\"64-57=\", \"62-1=\"], [\"3+70=\", \"70+12=\"], [\"46-15=\", \"69-68=\"], [\"17+38=\", \"74-16=\"], [\"51-46=\", \"31+47=\"], [\"98-50=\", \"51+28=\"], [\"77-2=\", \"16-6=\"], [\"75-57=\", \"95-54=\"], [\"31+13=\", \"64-0=\"], [\"44-22=\", \"1+85=\"], [\"19+9=\", \"15+66=\"], [\"8+12=\", \"74+5=\"], [\"95-35=\", \"68+26=\"], [\"7+17=\", \"91-71=\"], [\"22+36=\", \"54-50=\"], [\"5+56=\", \"69-8=\"], [\"35+47=\", \"95-47=\"], [\"86-35=\", \"62-43=\"], [\"85-28=\", \"43-34=\"], [\"73-39=\", \"56+2=\"], [\"93-80=\", \"51-35=\"], [\"54-24=\", \"98-66=\"], [\"60+11=\", \"18+62=\"], [\"44+24=\", \"38+30=\"], [\"8+66=\", \"81-23=\"], [\"60+11=\", \"77-15=\"], [\"92-33=\", \"95-85=\"], [\"14-2=\", \"35+21=\"], [\"76-64=\", \"28+28=\"], [\"16+33=\", \"50+32=\"], [\"26+72=\", \"4+36=\"], [\"58+15=\", \"8+63=\"], [\"53+42=\", \"22-12=\"], [\"82-79=\", \"18+47=\"], [\"30-21=\", \"32-2=\"], [\"86-58=\", \"21+48=\"], [\"70-20=\", \"41+19=\"], [\"24+33=\", \"88-37=\"], [\"94+1=\", \"58-31=\"], [\"3+41=\", \"70+10=\"], [\"9+46=\", \"57+31=\"], [\"71+6=\", \"79-8=\"], [\"59-39=\", \"55+38=\"], [\"71+12=\", \"93-25=\"], [\"76-54=\", \"6+61=\"], [\"47-42=\", \"23+8=\"], [\"38-23=\", \"14+64=\"], [\"50+8=\", \"37+6=\"], [\"98-19=\", \"4+54=\"], [\"22+20=\", \"9-9=\"], [\"57+38=\", \"7+87=\"], [\"85-28=\", \"97-20=\"], [\"29+17=\", \"13+2=\"], [\"93-47=\", \"88-54=\"], [\"84-13=\", \"8-3=\"], [\"49-43=\", \"39+15=\"], [\"92-6=\", \"85-60=\"]];\n\n// --- 1. The date paragraph is the first paragraph in the body. ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nconst [dateOld, dateNew] = replacements[0];\nif (dateParagraph.text.trim() === dateOld) {\n  dateParagraph.getRange().insertText(dateNew, Word.InsertLocation.replace);\n} else {\n  // Fall back to a document-wide search/replace if the structure differs\n  // from what we expect.\n  const results = context.document.body.search(dateOld, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(dateNew, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// --- 2. The arithmetic problems live in the lone table, 20 rows x 5 cols. ---\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst cellPairs = replacements.slice(1);\nconst columnCount = 5;\nconst rowCount = table.rowCount;\n\n// `cell.value` is directly settable and preserves the existing run and\n// paragraph formatting (font, size, alignment), unlike replacing via\n// insertText on the cell body, which would strip the cell's rPr/pPr.\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    if (idx >= cellPairs.length) break;\n    const [, newText] = cellPairs[idx];\n    const cell = table.getCell(r, c);\n    cell.value = newText;\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the date line and every \"a OP b=\" arithmetic-problem cell in the\n# single table with the values from the target revision. Replacements are\n# applied strictly by document position (first paragraph, then table cells\n# in row-major order) since several cell strings are not unique in the\n# source (e.g. \"85-28=\" and \"60+11=\" each occur twice), so plain\n# find/replace-by-text would be ambiguous.\n#\n# Each element is @(oldText, newText); element 0 is the date paragraph,\n# elements 1..100 are the 20x5 table cells, read left-to-right then\n# top-to-bottom.\n$replacements = @(\n    @('2023-08-12 Saturday', '2023-08-13 Sunday'),\n    @('13+0=', '0+26='),\n    @('16-11=', '93-65='),\n    @('93-4=', '60-12='),\n    @('26+9=', '92-43='),\n    @('93-71=', '74-22='),\n    @('82+2=', '74+12='),\n    @('58+3=', '21+75='),\n    @('50-14=', '0+40='),\n    @('66+17=', '11+31='),\n    @('13+18=', '76-46='),\n    @('9+40=', '68+31='),\n    @('75+8=', '66+8='),\n    @('83+1=', '18+20='),\n    @('38-29=', '14+81='),\n    @('56-15=', '32+24='),\n    @('0+90=', '52+18='),\n    @('11+66=', '37+37='),\n    @('87-19=', '52-49='),\n    @('80-17=', '70+25='),\n    @('20+34=', '30+47='),\n    @('77+12=', '47-29='),\n    @('57-37=', '9+75='),\n    @('48+48=', '53+22='),\n    @('12+79=', '81-34='),\n    @('67-42=', '16+12='),\n    @('72-20=', '67-22='),\n    @('27-22=', '19+8='),\n    @('88-40=', '1+53='),\n    @('8+70=', '23+18='),\n    @('46-24=', '30+68='),\n    @('96-67=', '31+64='),\n    @('91-46=', '88-29='),\n    @('41+58=', '68-59='),\n    @('2+22=', '39+11='),\n    @('89+3=', '9+25='),\n    @('0+99=', '99-58='),\n    @('67-1=', '99-39='),\n    @('26+7=', '2+66='),\n    @('96+2=', '93+3='),\n    @('16-7=', '65-27='),\n    @('67-40=', '53-37='),\n    @('69-39=', '83-5='),\n    @('89-60=', '11-10='),\n    @('64-57=', '62-1='),\n    @('3+70=', '70+12='),\n    @('46-15=', '69-68='),\n    @('17+38=', '74-16='),\n    @('51-46=', '31+47='),\n    @('98-50=', '51+28='),\n    @('77-2=', '16-6='),\n    @('75-57=', '95-54='),\n    @('31+13=', '64-0='),\n    @('44-22=', '1+85='),\n    @('19+9=', '15+66='),\n    @('8+12=', '74+5='),\n    @('95-35=', '68+26='),\n    @('7+17=', '91-71='),\n    @('22+36=', '54-50='),\n    @('5+56=', '69-8='),\n    @('35+47=', '95-47='),\n    @('86-35=', '62-43='),\n    @('85-28=', '43-34='),\n    @('73-39=', '56+2='),\n    @('93-80=', '51-35='),\n    @('54-24=', '98-66='),\n    @('60+11=', '18+62='),\n    @('44+24=', '38+30='),\n    @('8+66=', '81-23='),\n    @('60+11=', '77-15='),\n    @('92-33=', '95-85='),\n    @('14-2=', '35+21='),\n    @('76-64=', '28+28='),\n    @('16+33=', '50+32='),\n    @('26+72=', '4+36='),\n    @('58+15=', '8+63='),\n    @('53+42=', '22-12='),\n    @('82-79=', '18+47='),\n    @('30-21=', '32-2='),\n    @('86-58=', '21+48='),\n    @('70-20=', '41+19='),\n    @('24+33=', '88-37='),\n    @('94+1=', '58-31='),\n    @('3+41=', '70+10='),\n    @('9+46=', '57+31='),\n    @('71+6=', '79-8='),\n    @('59-39=', '55+38='),\n    @('71+12=', '93-25='),\n    @('76-54=', '6+61='),\n    @('47-42=', '23+8='),\n    @('38-23=', '14+64='),\n    @('50+8=', '37+6='),\n    @('98-19=', '4+54='),\n    @('22+20=', '9-9='),\n    @('57+38=', '7+87='),\n    @('85-28=', '97-20='),\n    @('29+17=', '13+2='),\n    @('93-47=', '88-54='),\n    @('84-13=', '8-3='),\n    @('49-43=', '39+15='),\n    @('92-6=', '85-60=')\n)\n\n$d = $word.ActiveDocument\n\n# --- 1. The date paragraph is the first paragraph in the body. ---\n$dateOld = $replacements[0][0]\n$dateNew = $replacements[0][1]\n$p = $d.Paragraphs.Item(1)\n$pText = $p.Range.Text.TrimEnd([char]13, [char]7)\nif ($pText -eq $dateOld) {\n    $p.Range.Text = $dateNew\n} else {\n    # Fall back to a document-wide Find/Replace if the structure differs\n    # from what we expect (wdReplaceAll = 2).\n    $find = $d.Content.Find\n    $find.Text = $dateOld\n    $find.Replacement.Text = $dateNew\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n\n# --- 2. The arithmetic problems live in the lone table, 20 rows x 5 cols. ---\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$idx = 1\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($idx -gt ($replacements.Length - 1)) { continue }\n        $newText = $replacements[$idx][1]\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newText\n        $idx++\n    }\n}\n"}
